$wb = $excel.ActiveWorkbook

# Sheet "展览" — update view-count column F for several events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 83
$ws1.Cells.Item(3, 6).Value = 4009
$ws1.Cells.Item(10, 6).Value = 117
$ws1.Cells.Item(15, 6).Value = 2810
$ws1.Cells.Item(16, 6).Value = 193

# Sheet "全部类型" — same events duplicated here, same column F updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 83
$ws4.Cells.Item(3, 6).Value = 4009
$ws4.Cells.Item(11, 6).Value = 117
$ws4.Cells.Item(18, 6).Value = 2810
$ws4.Cells.Item(19, 6).Value = 193
